$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 79
$srcRow = 78

# Copy formatting (cell style) from the row above so the new row matches
# the existing sheet's look (e.g. the date style on column A).
$ws.Range("A$srcRow`:J$srcRow").Copy()
$ws.Range("A$row`:J$row").PasteSpecial(-4122)

$ws.Cells.Item($row, 1).Value = 45635
$ws.Cells.Item($row, 2).Value = 116.4121952
$ws.Cells.Item($row, 3).Value = 0.00170247
$ws.Cells.Item($row, 4).Value = 0.008850780000000001
$ws.Cells.Item($row, 5).Value = 0.06933635
$ws.Cells.Item($row, 6).Value = 12792.90181321
$ws.Cells.Item($row, 7).Value = 465.80531254
$ws.Cells.Item($row, 8).Value = 0.24
$ws.Cells.Item($row, 9).Value = 1.7904431
$ws.Cells.Item($row, 10).Value = 485.38834923
